$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N16").ClearContents()
$ws.Range("H16").Value = 1199
$ws.Range("I16").Value = 1199
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1199
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -969

$ws.Range("H116").Value = 3823.1667
$ws.Range("I116").Value = 1750.8334
$ws.Range("J116").Value = 4859.3335
$ws.Range("K116").Value = 1750.8334
$ws.Range("L116").Value = 4859.3335
$ws.Range("M116").Value = 1691.1666
$ws.Range("N116").Value = -11743.3335

$ws.Range("H138").Value = 2788.9575
$ws.Range("J138").Value = 3094.359
$ws.Range("L138").Value = 9283.076999999999
$ws.Range("N138").Value = -19563.077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N22").ClearContents()
$ws.Range("H22").Value = 1906.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0

$ws.Range("H32").Value = 25259.043
$ws.Range("I32").Value = 32111.412
$ws.Range("J32").Value = 5844
$ws.Range("K32").Value = 32111.412
$ws.Range("L32").Value = 5844
$ws.Range("M32").Value = -31824.412
$ws.Range("N32").Value = -6418

$ws.Range("H96").Value = 20344
$ws.Range("J96").Value = 20344
$ws.Range("L96").Value = 20344
$ws.Range("N96").Value = -25836

$ws.Range("H110").Value = 2420.6667
$ws.Range("I110").Value = 1683
$ws.Range("K110").Value = 1683
$ws.Range("M110").Value = 362

$ws.Range("H114").Value = 39700
$ws.Range("J114").Value = 39700
$ws.Range("L114").Value = 39700
$ws.Range("N114").Value = -48378

$ws.Range("H132").Value = 15069.41
$ws.Range("I132").Value = 2253.4062
$ws.Range("K132").Value = 6760.2186
$ws.Range("M132").Value = -4230.2186

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1341.9375
$ws.Range("I99").Value = 871.1111
$ws.Range("K99").Value = 871.1111
$ws.Range("M99").Value = 626.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5134.95
$ws.Range("I99").Value = 4399.875
$ws.Range("J99").Value = 5625
$ws.Range("K99").Value = 4399.875
$ws.Range("L99").Value = 5625
$ws.Range("M99").Value = -2901.875
$ws.Range("N99").Value = -8621

$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 3500
$ws.Range("K122").Value = 10500
$ws.Range("M122").Value = -8050

$ws.Range("H126").Value = 5134.95
$ws.Range("I126").Value = 4399.875
$ws.Range("J126").Value = 5625
$ws.Range("K126").Value = 13199.625
$ws.Range("L126").Value = 16875
$ws.Range("M126").Value = -10729.625
$ws.Range("N126").Value = -21815

$ws.Range("H132").Value = 21327.312
$ws.Range("I132").Value = 23912.926
$ws.Range("J132").Value = 7365
$ws.Range("K132").Value = 71738.77799999999
$ws.Range("L132").Value = 22095
$ws.Range("M132").Value = -69208.77799999999
$ws.Range("N132").Value = -27155

$ws.Range("H134").Value = 7673.6665
$ws.Range("I134").Value = 842.0833
$ws.Range("J134").Value = 35000
$ws.Range("K134").Value = 2526.2499
$ws.Range("L134").Value = 105000
$ws.Range("M134").Value = 8.750100000000202
$ws.Range("N134").Value = -110070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 4588.5835
$ws.Range("I94").Value = 633.6667
$ws.Range("J94").Value = 5906.8887
$ws.Range("K94").Value = 1901.0001
$ws.Range("L94").Value = 17720.6661
$ws.Range("M94").Value = -1225.0001
$ws.Range("N94").Value = -19072.6661

$ws.Range("H96").Value = 752257
$ws.Range("J96").Value = 752257
$ws.Range("L96").Value = 2256771
$ws.Range("N96").Value = -2260889

$ws.Range("H107").Value = 9291.385
$ws.Range("I107").Value = 33600
$ws.Range("J107").Value = 1998.8
$ws.Range("K107").Value = 100800
$ws.Range("L107").Value = 5996.4
$ws.Range("M107").Value = -98880
$ws.Range("N107").Value = -9836.4

$ws.Range("H113").Value = 13190.125
$ws.Range("I113").Value = 50350.5
$ws.Range("J113").Value = 803.3333
$ws.Range("K113").Value = 151051.5
$ws.Range("L113").Value = 2409.9999
$ws.Range("M113").Value = -148881.5
$ws.Range("N113").Value = -6749.9999

$ws.Range("H122").Value = 1284.55
$ws.Range("I122").Value = 388.6
$ws.Range("J122").Value = 1583.2
$ws.Range("K122").Value = 3497.4
$ws.Range("L122").Value = 14248.8
$ws.Range("M122").Value = -1047.4
$ws.Range("N122").Value = -19148.8

$ws.Range("H131").Value = 822.2
$ws.Range("I131").Value = 507.5
$ws.Range("J131").Value = 835.3125
$ws.Range("K131").Value = 1522.5
$ws.Range("L131").Value = 2505.9375
$ws.Range("M131").Value = 3517.5
$ws.Range("N131").Value = -12585.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9379.666999999999
$ws.Range("I80").Value = 15319.375
$ws.Range("J80").Value = 4627.9
$ws.Range("K80").Value = 15319.375
$ws.Range("L80").Value = 4627.9
$ws.Range("M80").Value = -14321.375
$ws.Range("N80").Value = -6623.9

$ws.Range("H83").Value = 9379.666999999999
$ws.Range("I83").Value = 15319.375
$ws.Range("J83").Value = 4627.9
$ws.Range("K83").Value = 76596.875
$ws.Range("L83").Value = 23139.5
$ws.Range("M83").Value = -71604.875
$ws.Range("N83").Value = -33123.5

$ws.Range("H107").Value = 952
$ws.Range("I107").Value = 356.375
$ws.Range("J107").Value = 1547.625
$ws.Range("K107").Value = 356.375
$ws.Range("L107").Value = 1547.625
$ws.Range("M107").Value = 1563.625
$ws.Range("N107").Value = -5387.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4726.8213
$ws.Range("I7").Value = 4897.864
$ws.Range("J7").Value = 4099.6665
$ws.Range("K7").Value = 4897.864
$ws.Range("L7").Value = 4099.6665
$ws.Range("M7").Value = -4785.864
$ws.Range("N7").Value = -4323.6665

$ws.Range("H22").Value = 1266.4546
$ws.Range("I22").Value = 1304.9
$ws.Range("J22").Value = 882
$ws.Range("K22").Value = 1304.9
$ws.Range("L22").Value = 882
$ws.Range("M22").Value = -1009.9
$ws.Range("N22").Value = -1472

$ws.Range("H27").Value = 1266.4546
$ws.Range("I27").Value = 1304.9
$ws.Range("J27").Value = 882
$ws.Range("K27").Value = 1304.9
$ws.Range("L27").Value = 882
$ws.Range("M27").Value = -1197.9
$ws.Range("N27").Value = -1096

$ws.Range("H46").Value = 1533.9131
$ws.Range("I46").Value = 1605.0555
$ws.Range("J46").Value = 1277.8
$ws.Range("K46").Value = 1605.0555
$ws.Range("L46").Value = 1277.8
$ws.Range("M46").Value = -1417.0555
$ws.Range("N46").Value = -1653.8

$ws.Range("H126").Value = 4726.8213
$ws.Range("I126").Value = 4897.864
$ws.Range("J126").Value = 4099.6665
$ws.Range("K126").Value = 14693.592
$ws.Range("L126").Value = 12298.9995
$ws.Range("M126").Value = -12223.592
$ws.Range("N126").Value = -17238.9995

$ws.Range("H136").Value = 2363.6667
$ws.Range("I136").Value = 1795.4
$ws.Range("J136").Value = 2931.9333
$ws.Range("K136").Value = 5386.200000000001
$ws.Range("L136").Value = 8795.7999
$ws.Range("M136").Value = -2836.200000000001
$ws.Range("N136").Value = -13895.7999
